# testrange.xlsx update — "add JSON and ClassTest to project"
#
# 1. The original "Sheet" had a red fill applied to A2:C2 (no values). That
#    formatting/content is removed entirely, shrinking the sheet back down
#    to a bare A1:A1 sheet.
# 2. A new worksheet named "Pi" is added (after "Sheet") holding the value
#    of Pi (3.14) in cell F5.
# 3. A worksheet-scoped defined name "Pi" (scoped to "Sheet", localSheetId 0)
#    is created, pointing at Pi!F5.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet")

# Remove the old red-filled (but empty) A2:C2 block entirely so the sheet
# collapses back down to a single-cell A1:A1 used range.
$ws1.Range("A2:C2").Delete()

# Add the new "Pi" worksheet right after "Sheet".
$wsPi = $wb.Worksheets.Add($null, $ws1)
$wsPi.Name = "Pi"
$wsPi.Range("F5").Value = 3.14

# Create the defined name "Pi" local to the "Sheet" worksheet, pointing at
# Pi!F5.
$ws1.Names.Add("Pi", "=Pi!F5")

# Keep "Sheet" as the active tab (adding the new sheet makes it active by
# default).
$ws1.Activate()
